# Roll the 90-day GSC export window forward by one day:
#   - drop the oldest date row (2025-10-23)
#   - shift every remaining row up by one
#   - append a new row for 2026-01-21 (Non-HTTPS=0, HTTPS=25)
# The data lives on the worksheet literally named "Chart" (sheet1.xml),
# not "Table" (which only holds the Issue/Validation/Pages headers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Deleting row 2 (the 2025-10-23 row) shifts rows 3..91 up to 2..90,
# preserving each cell's existing type/shared-string id/style as Excel
# natively does on a row delete+shift (no risk of the text "2025-xx-xx"
# being reinterpreted as a date serial, which a plain .Value write would do).
$ws.Rows.Item(2).Delete()

# Row 91 is now free (used range shrank to A1:C90). Write the new last day.
# Prefix the date with an apostrophe so it is stored as text (matches the
# existing column A cells, which are all shared-string dates, not real
# Excel dates) instead of being auto-coerced into a date serial number.
$ws.Range("A91").Value = "'2026-01-21"

# That text assignment leaves a one-off "quote prefix" style on A91; copy
# the plain (style 0) formatting from the row above back onto it so the
# cell matches the rest of the column exactly.
$ws.Range("A90").Copy()
$ws.Range("A91").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B91").Value = 0
$ws.Range("C91").Value = 25
